$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "HardMode" column header
$ws.Range("C1").Value = "HardMode"

# Populate HardMode = TRUE for every data row (2 through 42)
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 3).Value = $true
}

# Update the view state: scrolled position and active selection
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("E34").Select()
